$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45008
$ws.Range("J2").Value = 1750
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 19000
$ws.Range("M2").Value = 18500
$ws.Range("P2").Value = 1850

$ws.Range("D3").Value = 44959
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 21000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 21500
$ws.Range("P3").Value = 2150

$ws.Range("D4").Value = 44265
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 1550

$ws.Range("D5").Value = 44406
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("P5").Value = 1450

$ws.Range("D6").Value = 44694
$ws.Range("K6").Value = 16000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 16500
$ws.Range("P6").Value = 1650

$ws.Range("D7").Value = 45001
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 18000
$ws.Range("M7").Value = 19000
$ws.Range("P7").Value = 1900

$ws.Range("D8").Value = 44904
$ws.Range("M8").Value = 14500
$ws.Range("P8").Value = 1450

$ws.Range("D9").Value = 44428
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 1550

$ws.Range("D10").Value = 44943
$ws.Range("K10").Value = 16000
$ws.Range("M10").Value = 16500
$ws.Range("P10").Value = 1650

$ws.Range("D11").Value = 44218
$ws.Range("J11").Value = 320
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 10500
$ws.Range("P11").Value = 1050

$ws.Range("D12").Value = 44547
$ws.Range("K12").Value = 19000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 19500
$ws.Range("P12").Value = 1950

$ws.Range("D13").Value = 44358
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("P13").Value = 1450

$ws.Range("D14").Value = 44377
$ws.Range("J14").Value = 650
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14538
$ws.Range("P14").Value = 1454

$ws.Range("D15").Value = 44679
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("P15").Value = 1950

$ws.Range("D16").Value = 44980
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15500
$ws.Range("P16").Value = 1550

$ws.Range("D17").Value = 44460
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 1550

$ws.Range("D18").Value = 44882
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15550
$ws.Range("P18").Value = 1555

$ws.Range("D19").Value = 44893
$ws.Range("J19").Value = 1400
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15571
$ws.Range("P19").Value = 1557

$ws.Range("D20").Value = 44263
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("P20").Value = 1550

$ws.Range("D21").Value = 44160
$ws.Range("J21").Value = 360
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 10500
$ws.Range("P21").Value = 1050

$ws.Range("D22").Value = 44988
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 16357
$ws.Range("P22").Value = 1636

$ws.Range("D23").Value = 44925
$ws.Range("J23").Value = 250
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14600
$ws.Range("P23").Value = 1460

$ws.Range("D24").Value = 44847
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16500
$ws.Range("P24").Value = 1650

$ws.Range("D25").Value = 44918
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12250
$ws.Range("P25").Value = 1225

$ws.Range("D26").Value = 44914
$ws.Range("J26").Value = 100

$ws.Range("D27").Value = 44860
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14500
$ws.Range("P27").Value = 1450

$ws.Range("D28").Value = 44291
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 13000
$ws.Range("L28").Value = 14000
$ws.Range("M28").Value = 13500
$ws.Range("P28").Value = 1350

$ws.Range("D29").Value = 45051
$ws.Range("J29").Value = 750
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16333
$ws.Range("P29").Value = 1633

$ws.Range("D30").Value = 44972
$ws.Range("J30").Value = 550
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15636
$ws.Range("P30").Value = 1564

$ws.Range("D31").Value = 44580
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 18000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = 19000
$ws.Range("P31").Value = 1900

$ws.Range("D32").Value = 44441
$ws.Range("J32").Value = 300
$ws.Range("M32").Value = 15500
$ws.Range("P32").Value = 1550

$ws.Range("D33").Value = 44524
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 21000
$ws.Range("M33").Value = 20500
$ws.Range("P33").Value = 2050

$ws.Range("D34").Value = 44330
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 13000
$ws.Range("L34").Value = 14000
$ws.Range("M34").Value = 13500
$ws.Range("P34").Value = 1350

$ws.Range("D35").Value = 44890
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 16000
$ws.Range("L35").Value = 17000
$ws.Range("M35").Value = 16500
$ws.Range("P35").Value = 1650

$ws.Range("D36").Value = 44777
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 24000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 24500
$ws.Range("P36").Value = 2450

$ws.Range("D37").Value = 45076
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = 16000
$ws.Range("P37").Value = 1600

$ws.Range("D38").Value = 44727
$ws.Range("K38").Value = 18000
$ws.Range("L38").Value = 19000
$ws.Range("M38").Value = 18500
$ws.Range("P38").Value = 1850

$ws.Range("D39").Value = 44714
$ws.Range("J39").Value = 400
$ws.Range("K39").Value = 19000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = 19500
$ws.Range("P39").Value = 1950

$ws.Range("D40").Value = 44204
$ws.Range("J40").Value = 400

$ws.Range("D41").Value = 44644
$ws.Range("J41").Value = 300
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 21000
$ws.Range("M41").Value = 20500
$ws.Range("P41").Value = 2050

$ws.Range("D42").Value = 44942
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 14000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = 14500
$ws.Range("P42").Value = 1450
